# Applies the LinuxForHealth rebrand / IG-republish edits described by the
# commit "Deploying to gh-pages from @ LinuxForHealth/alvearie-fhir-ig@..."
# to StructureDefinition-sent-to-recipient.xlsx.
#
# The workbook has two sheets:
#   "Metadata" - simple Property/Value table
#   "Elements" - wide FHIR element grid

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Metadata sheet: rebrand ibm.com -> linuxforhealth.org, bump version,
# refresh publish date, and rename the publisher.
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/sent-to-recipient"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---------------------------------------------------------------------
# Elements sheet: rebrand the extension URLs that appear inline in the
# Type(s)/Fixed Value columns, and clear the stale root-level
# Constraint(s) entry (ele-1/ext-1 now only applies further down the
# element tree).
# ---------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

$elements.Range("AI2").Value = ""
$elements.Range("J5").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/sent-time}`n"
$elements.Range("J6").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/outcome}`n"
$elements.Range("Q7").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/sent-to-recipient"
